# Scheduled-runner data refresh: update cached market-board figures
# (currentAveragePrice* / LevePrice* / LeveProfit*) across all eight
# job sheets, matching the latest Universalis pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1600
$ws.Range("I2").Value = 650
$ws.Range("J2").Value = 3500
$ws.Range("K2").Value = 650
$ws.Range("L2").Value = 3500
$ws.Range("M2").Value = -537
$ws.Range("N2").Value = -3726

$ws.Range("H41").Value = 1138.2
$ws.Range("I41").Value = 845.75
$ws.Range("J41").Value = 1333.1666
$ws.Range("K41").Value = 845.75
$ws.Range("L41").Value = 1333.1666
$ws.Range("M41").Value = -405.75
$ws.Range("N41").Value = -2213.1666

$ws.Range("H70").Value = 17899.75
$ws.Range("J70").Value = 20171.143
$ws.Range("L70").Value = 60513.429
$ws.Range("N70").Value = -61053.429

$ws.Range("H73").Value = 17899.75
$ws.Range("J73").Value = 20171.143
$ws.Range("L73").Value = 60513.429
$ws.Range("N73").Value = -62385.429

$ws.Range("H98").Value = 830.2121
$ws.Range("I98").Value = 536.52
$ws.Range("K98").Value = 536.52
$ws.Range("M98").Value = 961.48

$ws.Range("H106").Value = 5500.75
$ws.Range("I106").Value = 5500.75
$ws.Range("K106").Value = 5500.75
$ws.Range("M106").Value = -4869.75

$ws.Range("H112").Value = 1867.1428
$ws.Range("I112").Value = 500
$ws.Range("J112").Value = 1907.3529
$ws.Range("K112").Value = 1500
$ws.Range("L112").Value = 5722.0587
$ws.Range("M112").Value = -392
$ws.Range("N112").Value = -7938.0587

$ws.Range("H122").Value = 830.2121
$ws.Range("I122").Value = 536.52
$ws.Range("K122").Value = 1609.56
$ws.Range("M122").Value = 840.4400000000001

$ws.Range("H132").Value = 1455.7778
$ws.Range("I132").Value = 1510.4445
$ws.Range("K132").Value = 4531.333500000001
$ws.Range("M132").Value = -2001.333500000001

$ws.Range("H135").Value = 45455560
$ws.Range("J135").Value = 166667630
$ws.Range("L135").Value = 1500008670
$ws.Range("N135").Value = -1500013740

$ws.Range("H137").Value = 2742.7896
$ws.Range("I137").Value = 2160.8
$ws.Range("J137").Value = 2950.6428
$ws.Range("K137").Value = 6482.400000000001
$ws.Range("L137").Value = 8851.928400000001
$ws.Range("M137").Value = -3932.400000000001
$ws.Range("N137").Value = -13951.9284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5147.6123
$ws.Range("I32").Value = 3968.575
$ws.Range("J32").Value = 10387.777
$ws.Range("K32").Value = 3968.575
$ws.Range("L32").Value = 10387.777
$ws.Range("M32").Value = -3681.575
$ws.Range("N32").Value = -10961.777

$ws.Range("H45").Value = 1482.8334
$ws.Range("I45").Value = 1100
$ws.Range("K45").Value = 1100
$ws.Range("M45").Value = -723

$ws.Range("H61").Value = 2232.75
$ws.Range("I61").Value = 2265.5
$ws.Range("K61").Value = 2265.5
$ws.Range("M61").Value = -2053.5

$ws.Range("H102").Value = 1989
$ws.Range("I102").Value = 1989
$ws.Range("K102").Value = 1989
$ws.Range("M102").Value = -367

$ws.Range("H110").Value = 142
$ws.Range("I110").Value = 142
$ws.Range("K110").Value = 142
$ws.Range("M110").Value = 1903

$ws.Range("H132").Value = 1481.5714
$ws.Range("I132").Value = 1449.3334
$ws.Range("J132").Value = 1675
$ws.Range("K132").Value = 4348.0002
$ws.Range("L132").Value = 5025
$ws.Range("M132").Value = -1818.0002
$ws.Range("N132").Value = -10085

$ws.Range("H136").Value = 2232.75
$ws.Range("I136").Value = 2265.5
$ws.Range("K136").Value = 6796.5
$ws.Range("M136").Value = -4246.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2948.2727
$ws.Range("I20").Value = 2775.8572
$ws.Range("K20").Value = 2775.8572
$ws.Range("M20").Value = -2528.8572

$ws.Range("H94").Value = 844
$ws.Range("I94").Value = 835.7778
$ws.Range("J94").Value = 868.6667
$ws.Range("K94").Value = 835.7778
$ws.Range("L94").Value = 868.6667
$ws.Range("M94").Value = -384.7778
$ws.Range("N94").Value = -1770.6667

$ws.Range("H99").Value = 750.6667
$ws.Range("I99").Value = 641.3333
$ws.Range("K99").Value = 641.3333
$ws.Range("M99").Value = 856.6667

$ws.Range("H107").Value = 994.5
$ws.Range("I107").Value = 941.75
$ws.Range("K107").Value = 941.75
$ws.Range("M107").Value = 978.25

$ws.Range("H134").Value = 8886.117
$ws.Range("I134").Value = 9683.666999999999
$ws.Range("K134").Value = 29051.001
$ws.Range("M134").Value = -26516.001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1305.8
$ws.Range("I22").Value = 433.33334
$ws.Range("K22").Value = 433.33334
$ws.Range("M22").Value = -83.33334000000002

$ws.Range("H31").Value = 2803.625
$ws.Range("I31").Value = 2651.6667
$ws.Range("J31").Value = 2999
$ws.Range("K31").Value = 2651.6667
$ws.Range("L31").Value = 2999
$ws.Range("M31").Value = -2356.6667
$ws.Range("N31").Value = -3589

$ws.Range("H34").Value = 2803.625
$ws.Range("I34").Value = 2651.6667
$ws.Range("J34").Value = 2999
$ws.Range("K34").Value = 2651.6667
$ws.Range("L34").Value = 2999
$ws.Range("M34").Value = -2449.6667
$ws.Range("N34").Value = -3403

$ws.Range("H94").Value = 1900.5
$ws.Range("I94").Value = 1404
$ws.Range("J94").Value = 2397
$ws.Range("K94").Value = 1404
$ws.Range("L94").Value = 2397
$ws.Range("M94").Value = -953
$ws.Range("N94").Value = -3299

$ws.Range("H107").Value = 742.94116
$ws.Range("I107").Value = 499.06668
$ws.Range("J107").Value = 2572
$ws.Range("K107").Value = 499.06668
$ws.Range("L107").Value = 2572
$ws.Range("M107").Value = 1420.93332
$ws.Range("N107").Value = -6412

$ws.Range("H132").Value = 3065.375
$ws.Range("I132").Value = 1302
$ws.Range("K132").Value = 3906
$ws.Range("M132").Value = -1376

$ws.Range("H134").Value = 3705
$ws.Range("I134").Value = 3007.5715
$ws.Range("J134").Value = 5332.3335
$ws.Range("K134").Value = 9022.7145
$ws.Range("L134").Value = 15997.0005
$ws.Range("M134").Value = -6487.7145
$ws.Range("N134").Value = -21067.0005

$ws.Range("H141").Value = 71665
$ws.Range("J141").Value = 71497.5
$ws.Range("L141").Value = 71497.5
$ws.Range("N141").Value = -81857.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 999
$ws.Range("J45").Value = 999
$ws.Range("L45").Value = 2997
$ws.Range("N45").Value = -4061

$ws.Range("H131").Value = 12560.25
$ws.Range("I131").Value = 725
$ws.Range("J131").Value = 13405.625
$ws.Range("K131").Value = 2175
$ws.Range("L131").Value = 40216.875
$ws.Range("M131").Value = 2865
$ws.Range("N131").Value = -50296.875

$ws.Range("H132").Value = 1725.125
$ws.Range("I132").Value = 935.6667
$ws.Range("J132").Value = 2198.8
$ws.Range("K132").Value = 8421.0003
$ws.Range("L132").Value = 19789.2
$ws.Range("M132").Value = -5891.0003
$ws.Range("N132").Value = -24849.2

$ws.Range("H137").Value = 4551.8
$ws.Range("I137").Value = 1352
$ws.Range("J137").Value = 7169.8184
$ws.Range("K137").Value = 4056
$ws.Range("L137").Value = 21509.4552
$ws.Range("M137").Value = 1044
$ws.Range("N137").Value = -31709.4552

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2220
$ws.Range("I97").Value = 2090
$ws.Range("K97").Value = 2090
$ws.Range("M97").Value = -1594

$ws.Range("H102").Value = 7808.375
$ws.Range("I102").Value = 10395.6
$ws.Range("J102").Value = 3496.3333
$ws.Range("K102").Value = 10395.6
$ws.Range("L102").Value = 3496.3333
$ws.Range("M102").Value = -8773.6
$ws.Range("N102").Value = -6740.3333

$ws.Range("H122").Value = 2136.1428
$ws.Range("I122").Value = 2118
$ws.Range("K122").Value = 6354
$ws.Range("M122").Value = -3904

$ws.Range("H132").Value = 1482788.4
$ws.Range("I132").Value = 1926130
$ws.Range("J132").Value = 4982.8335
$ws.Range("K132").Value = 5778390
$ws.Range("L132").Value = 14948.5005
$ws.Range("M132").Value = -5775860
$ws.Range("N132").Value = -20008.5005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2974.9092
$ws.Range("I7").Value = 2764.1428
$ws.Range("K7").Value = 2764.1428
$ws.Range("M7").Value = -2652.1428

$ws.Range("H46").Value = 1461.9445
$ws.Range("I46").Value = 1096.125
$ws.Range("J46").Value = 1754.6
$ws.Range("K46").Value = 1096.125
$ws.Range("L46").Value = 1754.6
$ws.Range("M46").Value = -908.125
$ws.Range("N46").Value = -2130.6

$ws.Range("H82").Value = 1794.091
$ws.Range("I82").Value = 1307.2222
$ws.Range("K82").Value = 1307.2222
$ws.Range("M82").Value = -946.2221999999999

$ws.Range("H85").Value = 1794.091
$ws.Range("I85").Value = 1307.2222
$ws.Range("K85").Value = 1307.2222
$ws.Range("M85").Value = -59.22219999999993

$ws.Range("H100").Value = 3091.8
$ws.Range("I100").Value = 1833
$ws.Range("K100").Value = 1833
$ws.Range("M100").Value = -1292

$ws.Range("H126").Value = 2974.9092
$ws.Range("I126").Value = 2764.1428
$ws.Range("K126").Value = 8292.428400000001
$ws.Range("M126").Value = -5822.428400000001

$ws.Range("H132").Value = 2536.1875
$ws.Range("I132").Value = 1298
$ws.Range("J132").Value = 4128.143
$ws.Range("K132").Value = 3894
$ws.Range("L132").Value = 12384.429
$ws.Range("M132").Value = -1364
$ws.Range("N132").Value = -17444.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1914.5
$ws.Range("I81").Value = 699
$ws.Range("K81").Value = 1398
$ws.Range("M81").Value = -337

$ws.Range("H84").Value = 1914.5
$ws.Range("I84").Value = 699
$ws.Range("K84").Value = 6990
$ws.Range("M84").Value = -1686

$ws.Range("H132").Value = 1896.08
$ws.Range("I132").Value = 1238.421
$ws.Range("K132").Value = 3715.263
$ws.Range("M132").Value = -1185.263
